# Updates cryptos list values per the commit diff (prices + volume deltas,
# plus a few rows whose coin/link/price/volume shifted position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are plain text already, so a direct
# Value assignment is safe for them.
# Column D (Price) holds numeric-looking text (e.g. "59.004.87", "1.80")
# that must stay text -- Excel would otherwise coerce it to a Number and
# silently drop formatting like trailing zeros. Force text storage with a
# "@" number format, then restore the Normal style so no stray formatting
# is left behind.
# Column E (Volume) values already contain spaces/% so Excel keeps them text.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.004.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -6.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.424.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -9.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.566'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0985'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.07%  '
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.23'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.852.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.928.35'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.27%  '
$ws.Range('E16').Value = '  -7.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.478.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -7.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '321.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.964'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.47%  '
$ws.Range('E22').Value = '  -10.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.461'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '59.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.22%  '
$ws.Range('E25').Value = '  -5.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.972'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.52%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.91%  '
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0758'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -12.20%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.77%  '
$ws.Range('E37').Value = '  -3.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '306.75'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.826'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -14.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0933'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.573'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0521'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.77%  '
$ws.Range('E48').Value = '  -6.00%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.34%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -11.07%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.948.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.30%  '
